$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 (La Ligua, date 44890) swap with rows 6-9 (Provincia de Limarí, date 44908).
# Columns affected: D (Fecha), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado), R (Origen), S (Precio $/Kg).

$pairs = @(
    @{ a = 2; b = 6 },
    @{ a = 3; b = 7 },
    @{ a = 4; b = 8 },
    @{ a = 5; b = 9 }
)

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($pair in $pairs) {
    $rowA = $pair.a
    $rowB = $pair.b

    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}
